$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helpers
#
# Several cells in the "Crime Complaints" grid flip between the numeric-count
# style (s=15 / "#,##0"), the percent-style (s=16 / "#,##0.0") and the
# placeholder TEXT style (s=14 / General, holding the shared strings "0" or
# "***.*") as the weekly figures are refreshed. Simply assigning .Value keeps
# whatever xf the cell already had, which is correct for same-type edits but
# wrong when the cell's fundamental type changes (number <-> text). For those
# we copy over the number format from a cell that already has the desired
# style, using Copy/PasteSpecial(xlPasteFormats) so the workbook reuses the
# existing style record instead of inventing a new one.
# ---------------------------------------------------------------------------

function Set-NumberCell($targetAddr, $value, $refAddr) {
    # Target currently holds TEXT ("0" / "***.*"); turn it into a plain number,
    # copying number formatting from $refAddr (already in the desired style).
    $ws.Range($targetAddr).Value = $value
    $ws.Range($refAddr).Copy() | Out-Null
    $ws.Range($targetAddr).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
}

function Set-TextCell($targetAddr, $text, $refAddr) {
    # Target currently holds a NUMBER; turn it into the placeholder TEXT value,
    # copying number formatting from $refAddr (already in the desired style).
    $ws.Range($targetAddr).NumberFormat = "@"
    $ws.Range($targetAddr).Value = $text
    $ws.Range($refAddr).Copy() | Out-Null
    $ws.Range($targetAddr).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
}

# ---------------------------------------------------------------------------
# Masthead text tweaks (new report week)
# ---------------------------------------------------------------------------

# "Volume 31   Number  17" -> "...18"
$ws.Range("A8").Characters(21, 2).Text = "18"

# "Report Covering the Week  4/22/2024  Through  4/28/2024"
#                             ^ idx 27, len 9            ^ idx 47, len 9
$ws.Range("C9").Characters(27, 9).Text = "4/29/2024"
$ws.Range("C9").Characters(47, 9).Text = "5/5/2024"

# ---------------------------------------------------------------------------
# Row 14 - Murder
# ---------------------------------------------------------------------------
Set-NumberCell "C14" 1 "J14"
Set-NumberCell "F14" 1 "J14"
Set-NumberCell "I14" 1 "J14"
$ws.Range("K14").Value = -50
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -50
$ws.Range("N14").Value = -75

# ---------------------------------------------------------------------------
# Row 15 - Rape
# ---------------------------------------------------------------------------
Set-TextCell "C15" "0" "D15"

# ---------------------------------------------------------------------------
# Row 16 - Robbery
# ---------------------------------------------------------------------------
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 100
$ws.Range("F16").Value = 17
$ws.Range("G16").Value = 13
$ws.Range("H16").Value = 30.769230769230
$ws.Range("I16").Value = 65
$ws.Range("J16").Value = 48
$ws.Range("K16").Value = 35.416666666666
$ws.Range("L16").Value = 1.5625
$ws.Range("M16").Value = -27.777777777777
$ws.Range("N16").Value = -83.830845771144

# ---------------------------------------------------------------------------
# Row 17 - Fel. Assault
# ---------------------------------------------------------------------------
$ws.Range("C17").Value = 10
$ws.Range("D17").Value = 9
$ws.Range("E17").Value = 11.111111111111
$ws.Range("F17").Value = 32
$ws.Range("G17").Value = 32
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 111
$ws.Range("J17").Value = 135
$ws.Range("K17").Value = -17.777777777777
$ws.Range("L17").Value = 4.716981132075
$ws.Range("M17").Value = 56.338028169014
$ws.Range("N17").Value = 3.738317757009

# ---------------------------------------------------------------------------
# Row 18 - Burglary
# ---------------------------------------------------------------------------
Set-TextCell "C18" "0" "D15"
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -100
$ws.Range("F18").Value = 10
$ws.Range("G18").Value = 15
$ws.Range("H18").Value = -33.333333333333
$ws.Range("J18").Value = 50
$ws.Range("K18").Value = -30
$ws.Range("L18").Value = -16.666666666666
$ws.Range("M18").Value = -59.770114942528
$ws.Range("N18").Value = -92.341356673960

# ---------------------------------------------------------------------------
# Row 19 - Gr. Larceny
# ---------------------------------------------------------------------------
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 20
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 114
$ws.Range("J19").Value = 104
$ws.Range("K19").Value = 9.615384615384
$ws.Range("L19").Value = -11.627906976744
$ws.Range("M19").Value = 4.587155963302
$ws.Range("N19").Value = -35.955056179775

# ---------------------------------------------------------------------------
# Row 20 - G.L.A.
# ---------------------------------------------------------------------------
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 150
$ws.Range("G20").Value = 16
$ws.Range("H20").Value = 68.75
$ws.Range("I20").Value = 86
$ws.Range("J20").Value = 87
$ws.Range("K20").Value = -1.149425287356
$ws.Range("L20").Value = 7.5
$ws.Range("M20").Value = 4.878048780487
$ws.Range("N20").Value = -92.991035044824

# ---------------------------------------------------------------------------
# Row 21 - TOTAL
# ---------------------------------------------------------------------------
$ws.Range("C21").Value = 29
$ws.Range("D21").Value = 24
$ws.Range("E21").Value = 20.833333333333
$ws.Range("F21").Value = 112
$ws.Range("G21").Value = 96
$ws.Range("H21").Value = 16.666666666666
$ws.Range("I21").Value = 422
$ws.Range("J21").Value = 429
$ws.Range("K21").Value = -1.631701631701
$ws.Range("L21").Value = -1.170960187353
$ws.Range("M21").Value = -4.740406320541
$ws.Range("N21").Value = -82.306079664570

# ---------------------------------------------------------------------------
# Row 24 - Petit Larceny
# ---------------------------------------------------------------------------
$ws.Range("C24").Value = 29
$ws.Range("D24").Value = 35
$ws.Range("E24").Value = -17.142857142857
$ws.Range("F24").Value = 118
$ws.Range("G24").Value = 116
$ws.Range("H24").Value = 1.724137931034
$ws.Range("I24").Value = 492
$ws.Range("J24").Value = 513
$ws.Range("K24").Value = -4.093567251461
$ws.Range("L24").Value = 4.016913319238
$ws.Range("M24").Value = 87.786259541984

# ---------------------------------------------------------------------------
# Row 25 - Retail Theft
# ---------------------------------------------------------------------------
$ws.Range("C25").Value = 16
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = 300
$ws.Range("F25").Value = 53
$ws.Range("G25").Value = 30
$ws.Range("H25").Value = 76.666666666666
$ws.Range("I25").Value = 233
$ws.Range("J25").Value = 156
$ws.Range("K25").Value = 49.358974358974
$ws.Range("L25").Value = 26.630434782608

# ---------------------------------------------------------------------------
# Row 26 - Misd. Assault
# ---------------------------------------------------------------------------
$ws.Range("C26").Value = 19
$ws.Range("D26").Value = 7
$ws.Range("E26").Value = 171.428571428571
$ws.Range("F26").Value = 60
$ws.Range("G26").Value = 33
$ws.Range("H26").Value = 81.818181818181
$ws.Range("I26").Value = 220
$ws.Range("J26").Value = 176
$ws.Range("K26").Value = 25
$ws.Range("L26").Value = 34.146341463414
$ws.Range("M26").Value = 11.675126903553

# ---------------------------------------------------------------------------
# Row 27 - UCR Rape*
# ---------------------------------------------------------------------------
Set-TextCell "C27" "0" "D15"
Set-NumberCell "D27" 1 "J14"
Set-NumberCell "E27" -100 "K14"
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 150
$ws.Range("J27").Value = 9
$ws.Range("K27").Value = 33.333333333333

# ---------------------------------------------------------------------------
# Row 28 - Other Sex Crimes
# ---------------------------------------------------------------------------
Set-NumberCell "C28" 2 "J14"
Set-TextCell "D28" "0" "D15"
Set-TextCell "E28" "***.*" "E15"
$ws.Range("F28").Value = 4
$ws.Range("H28").Value = 33.333333333333
$ws.Range("I28").Value = 14
$ws.Range("K28").Value = -26.315789473684
$ws.Range("L28").Value = -6.666666666666

# ---------------------------------------------------------------------------
# Row 29 - Shooting Vic.
# ---------------------------------------------------------------------------
Set-NumberCell "C29" 5 "J14"
Set-NumberCell "F29" 5 "J14"
$ws.Range("I29").Value = 6
$ws.Range("K29").Value = 50
$ws.Range("L29").Value = -14.285714285714
$ws.Range("M29").Value = 100
$ws.Range("N29").Value = -57.142857142857

# ---------------------------------------------------------------------------
# Row 30 - Shooting Inc.
# ---------------------------------------------------------------------------
Set-NumberCell "C30" 2 "J14"
Set-NumberCell "F30" 2 "J14"
$ws.Range("I30").Value = 3
$ws.Range("K30").Value = -25
$ws.Range("L30").Value = -40
$ws.Range("M30").Value = 0
$ws.Range("N30").Value = -70
